$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7212023138999939
$ws.Range("B1").Value = 1.873716235160828
$ws.Range("C1").Value = 6.667419910430908
$ws.Range("D1").Value = 1.591243982315063
$ws.Range("E1").Value = 0.9138869047164917
